$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet index 1
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F3").Value = 43
$wsExpo.Range("F5").Value = 179
$wsExpo.Range("F6").Value = 1061
$wsExpo.Range("F7").Value = 1037
$wsExpo.Range("F8").Value = 8068
$wsExpo.Range("F9").Value = 131
$wsExpo.Range("F10").Value = 198
$wsExpo.Range("F11").Value = 6859
$wsExpo.Range("F14").Value = 4930
$wsExpo.Range("F17").Value = 5349
$wsExpo.Range("F18").Value = 1069
$wsExpo.Range("F19").Value = 323
$wsExpo.Range("F20").Value = 327
$wsExpo.Range("F21").Value = 447
$wsExpo.Range("F22").Value = 312
$wsExpo.Range("F26").Value = 96
$wsExpo.Range("F27").Value = 9066
$wsExpo.Range("F28").Value = 70
$wsExpo.Range("F29").Value = 1627
$wsExpo.Range("F30").Value = 68
$wsExpo.Range("F31").Value = 40
$wsExpo.Range("F33").Value = 838
$wsExpo.Range("F35").Value = 74
$wsExpo.Range("F37").Value = 1172
$wsExpo.Range("F39").Value = 4722
$wsExpo.Range("F40").Value = 30
$wsExpo.Range("F44").Value = 144
$wsExpo.Range("F45").Value = 72
$wsExpo.Range("F46").Value = 33
$wsExpo.Range("F47").Value = 1243
$wsExpo.Range("F48").Value = 35

# Sheet "演出" (Performances) - sheet index 2
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F9").Value = 178
$wsShow.Range("F11").Value = 6
$wsShow.Range("F14").Value = 86
$wsShow.Range("F17").Value = 888

# Sheet "全部类型" (All types) - sheet index 4
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value = 43
$wsAll.Range("F6").Value = 179
$wsAll.Range("F8").Value = 1061
$wsAll.Range("F9").Value = 1037
$wsAll.Range("F10").Value = 8068
$wsAll.Range("F11").Value = 131
$wsAll.Range("F12").Value = 198
$wsAll.Range("F13").Value = 6859
$wsAll.Range("F17").Value = 4930
$wsAll.Range("F19").Value = 5349
$wsAll.Range("F20").Value = 1069
$wsAll.Range("F21").Value = 323
$wsAll.Range("F22").Value = 327
$wsAll.Range("F23").Value = 447
$wsAll.Range("F24").Value = 312
$wsAll.Range("F28").Value = 96
$wsAll.Range("F29").Value = 178
$wsAll.Range("F30").Value = 9066
$wsAll.Range("F31").Value = 70
$wsAll.Range("F32").Value = 1627
$wsAll.Range("F33").Value = 40
$wsAll.Range("F35").Value = 838
$wsAll.Range("F37").Value = 74
$wsAll.Range("F39").Value = 1172
$wsAll.Range("F41").Value = 4722
$wsAll.Range("F44").Value = 144
$wsAll.Range("F45").Value = 72
$wsAll.Range("F46").Value = 33
$wsAll.Range("F47").Value = 1243
$wsAll.Range("F48").Value = 35

Write-Host "Update complete"
